# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the "全部类型" sheet, which carry duplicate data rows.
#
# Mapping of row -> new value (column F) per sheet:
#   展览 (sheet1):     F3=85  F4=106  F9=657->736  F12=10292  F16=646  F20=104
#   全部类型 (sheet4):  F3=85  F4=106  F10=736      F13=10292  F17=646  F21=104

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"    = @{ 3 = 85; 4 = 106; 9 = 736; 12 = 10292; 16 = 646; 20 = 104 }
    "全部类型" = @{ 3 = 85; 4 = 106; 10 = 736; 13 = 10292; 17 = 646; 21 = 104 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $newValue = $rowsMap[$row]
        $ws.Cells.Item($row, 6).Value = $newValue
    }
}
